# Auto-applied cell updates per commit diff (scheduled runner sync of Leve price/profit data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 59.333332   # H6: 51.714287 -> 59.333332
$ws.Cells.Item(6, 9).Value = 15.333333   # I6: 13 -> 15.333333
$ws.Cells.Item(6, 11).Value = 45.999999   # K6: 39 -> 45.999999
$ws.Cells.Item(6, 13).Value = 66.000001   # M6: 73 -> 66.000001
$ws.Cells.Item(34, 8).Value = 1944.4546   # H34: 2213.2727 -> 1944.4546
$ws.Cells.Item(34, 9).Value = 1944.4546   # I34: 2213.2727 -> 1944.4546
$ws.Cells.Item(34, 11).Value = 1944.4546   # K34: 2213.2727 -> 1944.4546
$ws.Cells.Item(34, 13).Value = -1741.4546   # M34: -2010.2727 -> -1741.4546
$ws.Cells.Item(36, 8).Value = 1944.4546   # H36: 2213.2727 -> 1944.4546
$ws.Cells.Item(36, 9).Value = 1944.4546   # I36: 2213.2727 -> 1944.4546
$ws.Cells.Item(36, 11).Value = 1944.4546   # K36: 2213.2727 -> 1944.4546
$ws.Cells.Item(36, 13).Value = -1229.4546   # M36: -1498.2727 -> -1229.4546
$ws.Cells.Item(41, 8).Value = 1616.3334   # H41: 230 -> 1616.3334
$ws.Cells.Item(41, 9).Value = 1225   # I41: 133.33333 -> 1225
$ws.Cells.Item(41, 10).Value = 2399   # J41: 375 -> 2399
$ws.Cells.Item(41, 11).Value = 1225   # K41: 133.33333 -> 1225
$ws.Cells.Item(41, 12).Value = 2399   # L41: 375 -> 2399
$ws.Cells.Item(41, 13).Value = -785   # M41: 306.66667 -> -785
$ws.Cells.Item(41, 14).Value = -3279   # N41: -1255 -> -3279
$ws.Cells.Item(55, 8).Value = 53.64706   # H55: 59.705883 -> 53.64706
$ws.Cells.Item(55, 9).Value = 78.28570999999999   # I55: 91.14286 -> 78.28570999999999
$ws.Cells.Item(55, 10).Value = 36.4   # J55: 37.7 -> 36.4
$ws.Cells.Item(55, 11).Value = 78.28570999999999   # K55: 91.14286 -> 78.28570999999999
$ws.Cells.Item(55, 12).Value = 36.4   # L55: 37.7 -> 36.4
$ws.Cells.Item(55, 13).Value = 135.71429   # M55: 122.85714 -> 135.71429
$ws.Cells.Item(55, 14).Value = -464.4   # N55: -465.7 -> -464.4
$ws.Cells.Item(98, 8).Value = 1876.5834   # H98: 1903.8462 -> 1876.5834
$ws.Cells.Item(98, 9).Value = 1910.2727   # I98: 1947.091 -> 1910.2727
$ws.Cells.Item(98, 10).Value = 1506   # J98: 1666 -> 1506
$ws.Cells.Item(98, 11).Value = 1910.2727   # K98: 1947.091 -> 1910.2727
$ws.Cells.Item(98, 12).Value = 1506   # L98: 1666 -> 1506
$ws.Cells.Item(98, 13).Value = -412.2727   # M98: -449.0909999999999 -> -412.2727
$ws.Cells.Item(98, 14).Value = -4502   # N98: -4662 -> -4502
$ws.Cells.Item(99, 8).Value = 907.75   # H99: 1001.7143 -> 907.75
$ws.Cells.Item(99, 9).Value = 537.4286   # I99: 585.3333 -> 537.4286
$ws.Cells.Item(99, 11).Value = 1612.2858   # K99: 1755.9999 -> 1612.2858
$ws.Cells.Item(99, 13).Value = -114.2857999999999   # M99: -257.9999 -> -114.2857999999999
$ws.Cells.Item(100, 8).Value = 3435.2727   # H100: 3435.3635 -> 3435.2727
$ws.Cells.Item(100, 9).Value = 424   # I100: 424.16666 -> 424
$ws.Cells.Item(100, 11).Value = 424   # K100: 424.16666 -> 424
$ws.Cells.Item(100, 13).Value = 117   # M100: 116.83334 -> 117
$ws.Cells.Item(107, 8).Value = 198.91304   # H107: 195.79167 -> 198.91304
$ws.Cells.Item(107, 9).Value = 122.61905   # I107: 122.681816 -> 122.61905
$ws.Cells.Item(107, 11).Value = 122.61905   # K107: 122.681816 -> 122.61905
$ws.Cells.Item(107, 13).Value = 1797.38095   # M107: 1797.318184 -> 1797.38095
$ws.Cells.Item(122, 8).Value = 1876.5834   # H122: 1903.8462 -> 1876.5834
$ws.Cells.Item(122, 9).Value = 1910.2727   # I122: 1947.091 -> 1910.2727
$ws.Cells.Item(122, 10).Value = 1506   # J122: 1666 -> 1506
$ws.Cells.Item(122, 11).Value = 5730.8181   # K122: 5841.272999999999 -> 5730.8181
$ws.Cells.Item(122, 12).Value = 4518   # L122: 4998 -> 4518
$ws.Cells.Item(122, 13).Value = -3280.8181   # M122: -3391.272999999999 -> -3280.8181
$ws.Cells.Item(122, 14).Value = -9418   # N122: -9898 -> -9418

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11024.564   # H32: 11851.694 -> 11024.564
$ws.Cells.Item(32, 9).Value = 6907.8486   # I32: 7488.7334 -> 6907.8486
$ws.Cells.Item(32, 11).Value = 6907.8486   # K32: 7488.7334 -> 6907.8486
$ws.Cells.Item(32, 13).Value = -6620.8486   # M32: -7201.7334 -> -6620.8486
$ws.Cells.Item(115, 8).Value = 54999   # H115: 0 -> 54999
$ws.Cells.Item(115, 10).Value = 54999   # J115: 0 -> 54999
$ws.Cells.Item(115, 12).Value = 54999   # L115: 0 -> 54999
$ws.Cells.Item(115, 14).Value = -58133   # N115: None -> -58133
$ws.Cells.Item(137, 8).Value = 69999   # H137: 0 -> 69999
$ws.Cells.Item(137, 10).Value = 69999   # J137: 0 -> 69999
$ws.Cells.Item(137, 12).Value = 69999   # L137: 0 -> 69999
$ws.Cells.Item(137, 14).Value = -80199   # N137: None -> -80199

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 204.61111   # H80: 204.77777 -> 204.61111
$ws.Cells.Item(80, 9).Value = 200.25   # I80: 184.6 -> 200.25
$ws.Cells.Item(80, 10).Value = 205.85715   # J80: 212.53847 -> 205.85715
$ws.Cells.Item(80, 11).Value = 200.25   # K80: 184.6 -> 200.25
$ws.Cells.Item(80, 12).Value = 205.85715   # L80: 212.53847 -> 205.85715
$ws.Cells.Item(80, 13).Value = 797.75   # M80: 813.4 -> 797.75
$ws.Cells.Item(80, 14).Value = -2201.85715   # N80: -2208.53847 -> -2201.85715
$ws.Cells.Item(83, 8).Value = 204.61111   # H83: 204.77777 -> 204.61111
$ws.Cells.Item(83, 9).Value = 200.25   # I83: 184.6 -> 200.25
$ws.Cells.Item(83, 10).Value = 205.85715   # J83: 212.53847 -> 205.85715
$ws.Cells.Item(83, 11).Value = 1001.25   # K83: 923 -> 1001.25
$ws.Cells.Item(83, 12).Value = 1029.28575   # L83: 1062.69235 -> 1029.28575
$ws.Cells.Item(83, 13).Value = 3990.75   # M83: 4069 -> 3990.75
$ws.Cells.Item(83, 14).Value = -11013.28575   # N83: -11046.69235 -> -11013.28575
$ws.Cells.Item(109, 8).Value = 59999   # H109: 0 -> 59999
$ws.Cells.Item(109, 10).Value = 59999   # J109: 0 -> 59999
$ws.Cells.Item(109, 12).Value = 59999   # L109: 0 -> 59999
$ws.Cells.Item(109, 14).Value = -62773   # N109: None -> -62773
$ws.Cells.Item(116, 8).Value = 64495   # H116: 99995 -> 64495
$ws.Cells.Item(116, 10).Value = 64495   # J116: 99995 -> 64495
$ws.Cells.Item(116, 12).Value = 64495   # L116: 99995 -> 64495
$ws.Cells.Item(116, 14).Value = -73673   # N116: -109173 -> -73673
$ws.Cells.Item(129, 8).Value = 64999   # H129: 0 -> 64999
$ws.Cells.Item(129, 10).Value = 64999   # J129: 0 -> 64999
$ws.Cells.Item(129, 12).Value = 64999   # L129: 0 -> 64999
$ws.Cells.Item(129, 14).Value = -74999   # N129: None -> -74999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2042.7646   # H22: 2241.6667 -> 2042.7646
$ws.Cells.Item(22, 10).Value = 2509.8   # J22: 2999.5 -> 2509.8
$ws.Cells.Item(22, 12).Value = 2509.8   # L22: 2999.5 -> 2509.8
$ws.Cells.Item(22, 14).Value = -3209.8   # N22: -3699.5 -> -3209.8
$ws.Cells.Item(31, 8).Value = 7740.6665   # H31: 8035.467 -> 7740.6665
$ws.Cells.Item(31, 9).Value = 7740.6665   # I31: 8035.467 -> 7740.6665
$ws.Cells.Item(31, 11).Value = 7740.6665   # K31: 8035.467 -> 7740.6665
$ws.Cells.Item(31, 13).Value = -7445.6665   # M31: -7740.467 -> -7445.6665
$ws.Cells.Item(34, 8).Value = 7740.6665   # H34: 8035.467 -> 7740.6665
$ws.Cells.Item(34, 9).Value = 7740.6665   # I34: 8035.467 -> 7740.6665
$ws.Cells.Item(34, 11).Value = 7740.6665   # K34: 8035.467 -> 7740.6665
$ws.Cells.Item(34, 13).Value = -7538.6665   # M34: -7833.467 -> -7538.6665
$ws.Cells.Item(44, 8).Value = 10071   # H44: 0 -> 10071
$ws.Cells.Item(44, 10).Value = 10071   # J44: 0 -> 10071
$ws.Cells.Item(44, 12).Value = 10071   # L44: 0 -> 10071
$ws.Cells.Item(44, 14).Value = -10955   # N44: None -> -10955
$ws.Cells.Item(51, 8).Value = 56062.5   # H51: 30825 -> 56062.5
$ws.Cells.Item(51, 9).Value = 0   # I51: 14000 -> 0
$ws.Cells.Item(51, 11).Value = 0   # K51: 14000 -> 0
$ws.Cells.Item(51, 13).ClearContents()   # M51: remove (was -13264)
$ws.Cells.Item(59, 8).Value = 32692.6   # H59: 57283.453 -> 32692.6
$ws.Cells.Item(59, 10).Value = 34091.332   # J59: 61001.4 -> 34091.332
$ws.Cells.Item(59, 12).Value = 34091.332   # L59: 61001.4 -> 34091.332
$ws.Cells.Item(59, 14).Value = -36381.332   # N59: -63291.4 -> -36381.332
$ws.Cells.Item(61, 8).Value = 56062.5   # H61: 30825 -> 56062.5
$ws.Cells.Item(61, 9).Value = 0   # I61: 14000 -> 0
$ws.Cells.Item(61, 11).Value = 0   # K61: 14000 -> 0
$ws.Cells.Item(61, 13).ClearContents()   # M61: remove (was -13652)
$ws.Cells.Item(98, 8).Value = 54999   # H98: 0 -> 54999
$ws.Cells.Item(98, 10).Value = 54999   # J98: 0 -> 54999
$ws.Cells.Item(98, 12).Value = 54999   # L98: 0 -> 54999
$ws.Cells.Item(98, 14).Value = -59491   # N98: None -> -59491
$ws.Cells.Item(107, 8).Value = 537.7778   # H107: 579.125 -> 537.7778
$ws.Cells.Item(107, 9).Value = 498.57144   # I107: 547.1667 -> 498.57144
$ws.Cells.Item(107, 11).Value = 498.57144   # K107: 547.1667 -> 498.57144
$ws.Cells.Item(107, 13).Value = 1421.42856   # M107: 1372.8333 -> 1421.42856
$ws.Cells.Item(114, 8).Value = 44999   # H114: 0 -> 44999
$ws.Cells.Item(114, 10).Value = 44999   # J114: 0 -> 44999
$ws.Cells.Item(114, 12).Value = 44999   # L114: 0 -> 44999
$ws.Cells.Item(114, 14).Value = -53677   # N114: None -> -53677

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 0   # H7: 112.25 -> 0
$ws.Cells.Item(7, 9).Value = 0   # I7: 126.333336 -> 0
$ws.Cells.Item(7, 10).Value = 0   # J7: 70 -> 0
$ws.Cells.Item(7, 11).Value = 0   # K7: 379.000008 -> 0
$ws.Cells.Item(7, 12).Value = 0   # L7: 210 -> 0
$ws.Cells.Item(7, 13).ClearContents()   # M7: remove (was -267.000008)
$ws.Cells.Item(7, 14).ClearContents()   # N7: remove (was -434)
$ws.Cells.Item(80, 8).Value = 4777.92   # H80: 4669.1924 -> 4777.92
$ws.Cells.Item(80, 9).Value = 4628.8423   # I80: 4628.8945 -> 4628.8423
$ws.Cells.Item(80, 10).Value = 5250   # J80: 4778.5713 -> 5250
$ws.Cells.Item(80, 11).Value = 13886.5269   # K80: 13886.6835 -> 13886.5269
$ws.Cells.Item(80, 12).Value = 15750   # L80: 14335.7139 -> 15750
$ws.Cells.Item(80, 13).Value = -12950.5269   # M80: -12950.6835 -> -12950.5269
$ws.Cells.Item(80, 14).Value = -17622   # N80: -16207.7139 -> -17622
$ws.Cells.Item(83, 8).Value = 4777.92   # H83: 4669.1924 -> 4777.92
$ws.Cells.Item(83, 9).Value = 4628.8423   # I83: 4628.8945 -> 4628.8423
$ws.Cells.Item(83, 10).Value = 5250   # J83: 4778.5713 -> 5250
$ws.Cells.Item(83, 11).Value = 41659.58070000001   # K83: 41660.0505 -> 41659.58070000001
$ws.Cells.Item(83, 12).Value = 47250   # L83: 43007.14169999999 -> 47250
$ws.Cells.Item(83, 13).Value = -36979.58070000001   # M83: -36980.0505 -> -36979.58070000001
$ws.Cells.Item(83, 14).Value = -56610   # N83: -52367.14169999999 -> -56610
$ws.Cells.Item(92, 8).Value = 2984.8572   # H92: 4949.25 -> 2984.8572
$ws.Cells.Item(92, 9).Value = 5047   # I92: 5838.8 -> 5047
$ws.Cells.Item(92, 10).Value = 2160   # J92: 3466.6667 -> 2160
$ws.Cells.Item(92, 11).Value = 15141   # K92: 17516.4 -> 15141
$ws.Cells.Item(92, 12).Value = 6480   # L92: 10400.0001 -> 6480
$ws.Cells.Item(92, 13).Value = -13893   # M92: -16268.4 -> -13893
$ws.Cells.Item(92, 14).Value = -8976   # N92: -12896.0001 -> -8976

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 34364.715   # H20: 21249.75 -> 34364.715
$ws.Cells.Item(20, 10).Value = 38425.668   # J20: 25000 -> 38425.668
$ws.Cells.Item(20, 12).Value = 38425.668   # L20: 25000 -> 38425.668
$ws.Cells.Item(20, 14).Value = -38915.668   # N20: -25490 -> -38915.668
$ws.Cells.Item(29, 8).Value = 9007.333000000001   # H29: 7003.5 -> 9007.333000000001
$ws.Cells.Item(29, 9).Value = 0   # I29: 4999 -> 0
$ws.Cells.Item(29, 10).Value = 9007.333000000001   # J29: 9008 -> 9007.333000000001
$ws.Cells.Item(29, 11).Value = 0   # K29: 4999 -> 0
$ws.Cells.Item(29, 12).Value = 9007.333000000001   # L29: 9008 -> 9007.333000000001
$ws.Cells.Item(29, 13).ClearContents()   # M29: remove (was -4709)
$ws.Cells.Item(29, 14).Value = -9587.333000000001   # N29: -9588 -> -9587.333000000001
$ws.Cells.Item(40, 8).Value = 12100.2   # H40: 14250.25 -> 12100.2
$ws.Cells.Item(40, 9).Value = 13375.25   # I40: 16667 -> 13375.25
$ws.Cells.Item(40, 11).Value = 13375.25   # K40: 16667 -> 13375.25
$ws.Cells.Item(40, 13).Value = -13224.25   # M40: -16516 -> -13224.25
$ws.Cells.Item(43, 8).Value = 1606.5555   # H43: 2195.9 -> 1606.5555
$ws.Cells.Item(63, 8).Value = 29000   # H63: 30000 -> 29000
$ws.Cells.Item(63, 10).Value = 29000   # J63: 30000 -> 29000
$ws.Cells.Item(63, 12).Value = 29000   # L63: 30000 -> 29000
$ws.Cells.Item(63, 14).Value = -30372   # N63: -31372 -> -30372
$ws.Cells.Item(66, 8).Value = 29000   # H66: 30000 -> 29000
$ws.Cells.Item(66, 10).Value = 29000   # J66: 30000 -> 29000
$ws.Cells.Item(66, 12).Value = 87000   # L66: 90000 -> 87000
$ws.Cells.Item(66, 14).Value = -93864   # N66: -96864 -> -93864
$ws.Cells.Item(97, 8).Value = 1126.1818   # H97: 1080.4 -> 1126.1818
$ws.Cells.Item(97, 9).Value = 925.375   # I97: 925.5 -> 925.375
$ws.Cells.Item(97, 10).Value = 1661.6666   # J97: 1257.4286 -> 1661.6666
$ws.Cells.Item(97, 11).Value = 925.375   # K97: 925.5 -> 925.375
$ws.Cells.Item(97, 12).Value = 1661.6666   # L97: 1257.4286 -> 1661.6666
$ws.Cells.Item(97, 13).Value = -429.375   # M97: -429.5 -> -429.375
$ws.Cells.Item(97, 14).Value = -2653.6666   # N97: -2249.4286 -> -2653.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 2500   # H4: 4750 -> 2500
$ws.Cells.Item(4, 10).Value = 0   # J4: 7000 -> 0
$ws.Cells.Item(4, 12).Value = 0   # L4: 7000 -> 0
$ws.Cells.Item(4, 14).ClearContents()   # N4: remove (was -7226)
$ws.Cells.Item(25, 8).Value = 7   # H25: 0 -> 7
$ws.Cells.Item(25, 9).Value = 7   # I25: 0 -> 7
$ws.Cells.Item(25, 11).Value = 7   # K25: 0 -> 7
$ws.Cells.Item(25, 13).Value = 223   # M25: None -> 223
$ws.Cells.Item(28, 8).Value = 2500   # H28: 4750 -> 2500
$ws.Cells.Item(28, 10).Value = 0   # J28: 7000 -> 0
$ws.Cells.Item(28, 12).Value = 0   # L28: 7000 -> 0
$ws.Cells.Item(28, 14).ClearContents()   # N28: remove (was -7464)
$ws.Cells.Item(37, 8).Value = 2500   # H37: 4750 -> 2500
$ws.Cells.Item(37, 10).Value = 0   # J37: 7000 -> 0
$ws.Cells.Item(37, 12).Value = 0   # L37: 7000 -> 0
$ws.Cells.Item(37, 14).ClearContents()   # N37: remove (was -7214)
$ws.Cells.Item(53, 8).Value = 6999   # H53: 0 -> 6999
$ws.Cells.Item(53, 10).Value = 6999   # J53: 0 -> 6999
$ws.Cells.Item(53, 12).Value = 6999   # L53: 0 -> 6999
$ws.Cells.Item(53, 14).Value = -8035   # N53: None -> -8035
$ws.Cells.Item(68, 8).Value = 9999.75   # H68: 9999.799999999999 -> 9999.75
$ws.Cells.Item(71, 8).Value = 9999.75   # H71: 9999.799999999999 -> 9999.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 37336.332   # H7: 53504.5 -> 37336.332
$ws.Cells.Item(7, 9).Value = 6002   # I7: 7004 -> 6002
$ws.Cells.Item(7, 11).Value = 6002   # K7: 7004 -> 6002
$ws.Cells.Item(7, 13).Value = -5889   # M7: -6891 -> -5889
$ws.Cells.Item(21, 8).Value = 1302998.2   # H21: 1203099.4 -> 1302998.2
$ws.Cells.Item(21, 9).Value = 1624998   # I21: 1500124.5 -> 1624998
$ws.Cells.Item(21, 11).Value = 1624998   # K21: 1500124.5 -> 1624998
$ws.Cells.Item(21, 13).Value = -1624763   # M21: -1499889.5 -> -1624763
$ws.Cells.Item(35, 8).Value = 1302998.2   # H35: 1203099.4 -> 1302998.2
$ws.Cells.Item(35, 9).Value = 1624998   # I35: 1500124.5 -> 1624998
$ws.Cells.Item(35, 11).Value = 1624998   # K35: 1500124.5 -> 1624998
$ws.Cells.Item(35, 13).Value = -1624708   # M35: -1499834.5 -> -1624708
$ws.Cells.Item(37, 8).Value = 24999   # H37: 0 -> 24999
$ws.Cells.Item(37, 10).Value = 24999   # J37: 0 -> 24999
$ws.Cells.Item(37, 12).Value = 24999   # L37: 0 -> 24999
$ws.Cells.Item(37, 14).Value = -25405   # N37: None -> -25405
$ws.Cells.Item(63, 8).Value = 33333   # H63: 35000 -> 33333
$ws.Cells.Item(63, 10).Value = 32499.5   # J63: 0 -> 32499.5
$ws.Cells.Item(63, 12).Value = 32499.5   # L63: 0 -> 32499.5
$ws.Cells.Item(63, 14).Value = -33747.5   # N63: None -> -33747.5
$ws.Cells.Item(64, 8).Value = 59967.5   # H64: 59969 -> 59967.5
$ws.Cells.Item(64, 10).Value = 59967.5   # J64: 59969 -> 59967.5
$ws.Cells.Item(64, 12).Value = 59967.5   # L64: 59969 -> 59967.5
$ws.Cells.Item(64, 14).Value = -60463.5   # N64: -60465 -> -60463.5
$ws.Cells.Item(66, 8).Value = 33333   # H66: 35000 -> 33333
$ws.Cells.Item(66, 10).Value = 32499.5   # J66: 0 -> 32499.5
$ws.Cells.Item(66, 12).Value = 97498.5   # L66: 0 -> 97498.5
$ws.Cells.Item(66, 14).Value = -103738.5   # N66: None -> -103738.5
$ws.Cells.Item(67, 8).Value = 59967.5   # H67: 59969 -> 59967.5
$ws.Cells.Item(67, 10).Value = 59967.5   # J67: 59969 -> 59967.5
$ws.Cells.Item(67, 12).Value = 59967.5   # L67: 59969 -> 59967.5
$ws.Cells.Item(67, 14).Value = -61683.5   # N67: -61685 -> -61683.5
